# pico-glitcher-v1.1-BOM: rectify erroneous comment + refreshed JLCPCB part numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 -> C8, 22uF Aluminum Electrolytic Capacitor
$ws.Range("C5").Value = "SMD"
$ws.Range("D5").Value = "C178580"

# Row 10 -> R12, R13, R14, R15 (4x 10kOhm 1206 Resistor Array)
$ws.Range("C10").Value = "0603x4"
$ws.Range("D10").Value = "C396839"

# Row 11 -> R16, R17, R18, R19 (4x 100Ohm 1206 Resistor Array)
$ws.Range("D11").Value = "C396838"

# Row 26 -> S_RST, Tactile Switch
$ws.Range("D26").Value = "C2845324"
